$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..10) {
    $cell = $ws.Cells.Item($r, 3)  # column C
    if ($cell.Value2 -eq 45174) {
        $cell.Value2 = 45175
    }
}
